# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
#
# The "Date" column (BF) on every data row (2-31) held the literal text
# "5-11-2011-12" (a mangled concatenation of the source file name). Correct
# it to the real ISO game date "2012-05-11", keeping it as plain text
# (not an Excel date serial, and without touching cell formatting).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "5-11-2011-12"
$newValue = "2012-05-11"

# A direct `.Value = "2012-05-11"` assignment gets auto-recognized by Excel
# as a date and converted to a date serial number (plus a new date number
# format). To keep it as literal text, stage the text as a formula result
# in a scratch cell, then paste-special just the resulting value into each
# target cell - this bypasses Excel's "looks like a date" autoconversion.
$scratch = $ws.Range("ZZ1")
$scratch.Formula = "=""" + $newValue + """"
$scratch.Copy()

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Range("BF$row")
    if ($cell.Value2 -eq $oldValue) {
        $cell.PasteSpecial(-4163)  # xlPasteValues
    }
}

$scratch.ClearContents()
$excel.CutCopyMode = $false
